$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: correct a duplicated/incorrect value for this employee's hours ---
$ws.Range("F5").Value = 250
# N5 (VENCIMENTO_DB date) corrected
$ws.Range("N5").Value = 44899

# --- New row 6: additional entry (avoids future duplicate-data confusion) ---
# Copy row 5's formatting down to row 6 first so number formats (dates, etc.)
# match the rest of the table.
$ws.Range("A5:P5").Copy()
$ws.Range("A6:P6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "DBA"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "MANHÃ"
$ws.Range("D6").Value = "JUNIOR"
$ws.Range("E6").Value = "ALIMENTACAO"
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = "JUNIOR"
$ws.Range("I6").Value = "SOUZA"
$ws.Range("J6").Value = "RUA JOSÉ"
$ws.Range("K6").Value = 999238821
$ws.Range("L6").Value = 2324981042
$ws.Range("M6").Formula = "=F6*50%"
$ws.Range("N6").Value = 44899
$ws.Range("O6").Value = 50
$ws.Range("P6").Value = 44905

# --- Update the active selection to match where the user ended up ---
$ws.Range("N8").Select() | Out-Null
